# Refresh the crypto price/volume table with the latest scrape.
# Every data cell in the sheet is stored as text (the source feed renders
# prices/percentages as formatted strings, e.g. "65.495.69" / "  -3.11%  "),
# so writes below force/preserve text storage for any new value that looks
# like a plain number -- otherwise Excel's automatic type detection would
# silently convert the cell to a Number (and can even rewrite the literal
# digits, e.g. "1.00" -> 1, "137.60" -> 137.6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "65.495.69"; ForceText = $false },
    @{ Cell = "E2"; Value = "  -3.11%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "3.494.46"; ForceText = $false },
    @{ Cell = "E3"; Value = "  -0.32%  "; ForceText = $false },
    @{ Cell = "D4"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E4"; Value = "  +0.04%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "554.94"; ForceText = $true },
    @{ Cell = "E5"; Value = "  +0.07%  "; ForceText = $false },
    @{ Cell = "D6"; Value = "178.85"; ForceText = $true },
    @{ Cell = "E6"; Value = "  -6.37%  "; ForceText = $false },
    @{ Cell = "D7"; Value = "0.637"; ForceText = $true },
    @{ Cell = "E7"; Value = "  +4.05%  "; ForceText = $false },
    @{ Cell = "E8"; Value = "  +0.05%  "; ForceText = $false },
    @{ Cell = "D9"; Value = "0.630"; ForceText = $true },
    @{ Cell = "E9"; Value = "  -1.60%  "; ForceText = $false },
    @{ Cell = "D10"; Value = "0.154"; ForceText = $true },
    @{ Cell = "E10"; Value = "  +2.70%  "; ForceText = $false },
    @{ Cell = "D11"; Value = "53.69"; ForceText = $true },
    @{ Cell = "E11"; Value = "  -6.31%  "; ForceText = $false },
    @{ Cell = "D12"; Value = "0.0000271"; ForceText = $true },
    @{ Cell = "E12"; Value = "  -1.77%  "; ForceText = $false },
    @{ Cell = "D13"; Value = "9.24"; ForceText = $true },
    @{ Cell = "E13"; Value = "  -2.73%  "; ForceText = $false },
    @{ Cell = "D14"; Value = "4.063.41"; ForceText = $false },
    @{ Cell = "E14"; Value = "  -0.03%  "; ForceText = $false },
    @{ Cell = "D15"; Value = "3.501.09"; ForceText = $false },
    @{ Cell = "E15"; Value = "  +0.04%  "; ForceText = $false },
    @{ Cell = "B16"; Value = "TRON"; ForceText = $false },
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; ForceText = $false },
    @{ Cell = "D16"; Value = "0.121"; ForceText = $true },
    @{ Cell = "E16"; Value = "  +0.10%  "; ForceText = $false },
    @{ Cell = "B17"; Value = "Chainlink"; ForceText = $false },
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; ForceText = $false },
    @{ Cell = "D17"; Value = "18.39"; ForceText = $true },
    @{ Cell = "E17"; Value = "  +0.17%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "12.12"; ForceText = $true },
    @{ Cell = "E18"; Value = "  +2.25%  "; ForceText = $false },
    @{ Cell = "D19"; Value = "65.548.86"; ForceText = $false },
    @{ Cell = "E19"; Value = "  -4.13%  "; ForceText = $false },
    @{ Cell = "D20"; Value = "0.995"; ForceText = $true },
    @{ Cell = "E20"; Value = "  -1.60%  "; ForceText = $false },
    @{ Cell = "D21"; Value = "413.69"; ForceText = $true },
    @{ Cell = "E21"; Value = "  +1.80%  "; ForceText = $false },
    @{ Cell = "D22"; Value = "4.05"; ForceText = $true },
    @{ Cell = "E22"; Value = "  +2.50%  "; ForceText = $false },
    @{ Cell = "D23"; Value = "85.92"; ForceText = $true },
    @{ Cell = "E23"; Value = "  +1.29%  "; ForceText = $false },
    @{ Cell = "D24"; Value = "4.11"; ForceText = $true },
    @{ Cell = "E24"; Value = "  -2.69%  "; ForceText = $false },
    @{ Cell = "D25"; Value = "12.75"; ForceText = $true },
    @{ Cell = "E25"; Value = "  +7.41%  "; ForceText = $false },
    @{ Cell = "D26"; Value = "10.79"; ForceText = $true },
    @{ Cell = "E26"; Value = "  -9.09%  "; ForceText = $false },
    @{ Cell = "D27"; Value = "2.85"; ForceText = $true },
    @{ Cell = "E27"; Value = "  -2.17%  "; ForceText = $false },
    @{ Cell = "D28"; Value = "6.04"; ForceText = $true },
    @{ Cell = "E28"; Value = "  -1.47%  "; ForceText = $false },
    @{ Cell = "D29"; Value = "9.02"; ForceText = $true },
    @{ Cell = "E29"; Value = "  +4.33%  "; ForceText = $false },
    @{ Cell = "D30"; Value = "30.25"; ForceText = $true },
    @{ Cell = "E30"; Value = "  -0.82%  "; ForceText = $false },
    @{ Cell = "D31"; Value = "6.47"; ForceText = $true },
    @{ Cell = "E31"; Value = "  -6.00%  "; ForceText = $false },
    @{ Cell = "D32"; Value = "608.01"; ForceText = $true },
    @{ Cell = "E32"; Value = "  -11.60%  "; ForceText = $false },
    @{ Cell = "D33"; Value = "11.64"; ForceText = $true },
    @{ Cell = "E33"; Value = "  -0.54%  "; ForceText = $false },
    @{ Cell = "E34"; Value = "  -1.03%  "; ForceText = $false },
    @{ Cell = "D35"; Value = "59.58"; ForceText = $true },
    @{ Cell = "D36"; Value = "0.147"; ForceText = $true },
    @{ Cell = "E36"; Value = "  +9.62%  "; ForceText = $false },
    @{ Cell = "E37"; Value = "  -0.10%  "; ForceText = $false },
    @{ Cell = "B38"; Value = "InjectiveProtocol"; ForceText = $false },
    @{ Cell = "C38"; Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; ForceText = $false },
    @{ Cell = "D38"; Value = "37.22"; ForceText = $true },
    @{ Cell = "E38"; Value = "  -4.74%  "; ForceText = $false },
    @{ Cell = "B39"; Value = "PEPE"; ForceText = $false },
    @{ Cell = "C39"; Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; ForceText = $false },
    @{ Cell = "D39"; Value = "0.0₃0789"; ForceText = $false },
    @{ Cell = "E39"; Value = "  -4.45%  "; ForceText = $false },
    @{ Cell = "D40"; Value = "3.365.50"; ForceText = $false },
    @{ Cell = "E40"; Value = "  +10.08%  "; ForceText = $false },
    @{ Cell = "D41"; Value = "0.379"; ForceText = $true },
    @{ Cell = "E41"; Value = "  -6.16%  "; ForceText = $false },
    @{ Cell = "D42"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E42"; Value = "  -0.08%  "; ForceText = $false },
    @{ Cell = "D43"; Value = "3.26"; ForceText = $true },
    @{ Cell = "E43"; Value = "  -3.39%  "; ForceText = $false },
    @{ Cell = "E44"; Value = "  -4.97%  "; ForceText = $false },
    @{ Cell = "D45"; Value = "2.53"; ForceText = $true },
    @{ Cell = "E45"; Value = "  -9.93%  "; ForceText = $false },
    @{ Cell = "E46"; Value = "  -1.57%  "; ForceText = $false },
    @{ Cell = "D47"; Value = "3.23"; ForceText = $true },
    @{ Cell = "E47"; Value = "  +0.07%  "; ForceText = $false },
    @{ Cell = "D48"; Value = "2.71"; ForceText = $true },
    @{ Cell = "E48"; Value = "  -2.60%  "; ForceText = $false },
    @{ Cell = "D49"; Value = "0.132"; ForceText = $true },
    @{ Cell = "E49"; Value = "  +1.52%  "; ForceText = $false },
    @{ Cell = "D50"; Value = "8.44"; ForceText = $true },
    @{ Cell = "E50"; Value = "  -5.98%  "; ForceText = $false },
    @{ Cell = "D51"; Value = "137.60"; ForceText = $true },
    @{ Cell = "E51"; Value = "  -2.12%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Stamp the cell as Text first so the numeral is kept verbatim,
        # then hand the style back to Normal so no stray format lingers.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
